# Apply the "Constants" sheet changes:
#  - Insert a new row at row 7 with "MaxExecutionAttemptsHigh" / 99999 / description
#  - Insert a new row at row 10 (after the previous insert shifts things down) with
#    "RetryIntervalLow" / 1 / the same duration description used by "RetryInterval"
#
# Net effect matches the target diff: two rows inserted into the Constants sheet,
# pushing every following row down by 2, dimension growing from A1:Z64 to A1:Z66
# (and the trailing filler rows growing from 1021 to 1023).

$wb = $excel.ActiveWorkbook

$constants = $wb.Worksheets.Item("Constants")
$settings  = $wb.Worksheets.Item("Settings")

# --- Insert row for "MaxExecutionAttemptsHigh" above the existing "MaxLockTimeout" row ---
$constants.Rows.Item(7).Insert()
$constants.Range("A7").Value = "MaxExecutionAttemptsHigh"
$constants.Range("B7").Value = 99999
$constants.Range("C7").Value = "Maximum number of execution attempts for a process step which by default is high."

# --- Insert row for "RetryIntervalLow" right after "RetryInterval" (now row 9) ---
$constants.Rows.Item(10).Insert()
$constants.Range("A10").Value = "RetryIntervalLow"
$constants.Range("B10").Value = 1
$constants.Range("C10").Value = "Duration, in seconds, between re-execution attempts"
# Match the styling of the neighboring "RetryInterval" description cell (C9)
$constants.Range("C9").Copy()
$constants.Range("C10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply the text since PasteSpecial(formats) on C10 should only touch formatting, but make sure
# the value wasn't affected by the paste special call above.
$constants.Range("C10").Value = "Duration, in seconds, between re-execution attempts"

# --- Selections, matching the final state of the workbook ---
$constants.Range("A10:XFD10").Select()
$settings.Range("A24").Select()

$wb.Save()
